# Update the "想去人数" (want-to-go count) column F for the events that
# changed in this data refresh. The same rows/values need updating on both
# the "展览" sheet and the "全部类型" sheet (they mirror the same data).

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 246
    3  = 273
    4  = 283
    5  = 825
    6  = 279
    7  = 6711
    10 = 119
    11 = 83
    12 = 38
    13 = 12
    16 = 222
    17 = 569
    18 = 63
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
